$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 630
$ws.Range("F3").Value = 152
$ws.Range("F4").Value = 219
$ws.Range("F5").Value = 887
$ws.Range("F6").Value = 96
$ws.Range("F8").Value = 23
$ws.Range("F9").Value = 5657
$ws.Range("F11").Value = 880
$ws.Range("F12").Value = 20
$ws.Range("F13").Value = 691
$ws.Range("F14").Value = 417
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 352
$ws.Range("F17").Value = 791
$ws.Range("F18").Value = 1972
$ws.Range("F19").Value = 84
$ws.Range("F20").Value = 713
$ws.Range("F21").Value = 354
$ws.Range("F22").Value = 149

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 160

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 160
$ws.Range("F3").Value = 630
$ws.Range("F4").Value = 152
$ws.Range("F5").Value = 219
$ws.Range("F6").Value = 887
$ws.Range("F7").Value = 96
$ws.Range("F10").Value = 27
$ws.Range("F11").Value = 23
$ws.Range("F12").Value = 5657
$ws.Range("F14").Value = 880
$ws.Range("F15").Value = 20
$ws.Range("F16").Value = 691
$ws.Range("F17").Value = 417
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 352
$ws.Range("F20").Value = 791
$ws.Range("F21").Value = 1972
$ws.Range("F22").Value = 84
$ws.Range("F23").Value = 713
$ws.Range("F24").Value = 354
$ws.Range("F25").Value = 149

